$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells for team record columns
$ws.Range("AC1").Value = "Wins"
$ws.Range("AD1").Value = "Losses"
$ws.Range("AE1").Value = "Ties"

# Match the header formatting (bold, centered, bordered) used by existing headers
$ws.Range("AA1").Copy()
$ws.Range("AC1:AE1").PasteSpecial(-4122)

# Fill in the team record values for every data row (2-45)
$ws.Range("AC2:AC45").Value = 82
$ws.Range("AD2:AD45").Value = 80
$ws.Range("AE2:AE45").Value = 0

Write-Host "Added Wins/Losses/Ties columns (AC:AE)"
